$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 403
$ws.Range("A403").Value = '2025-11-28T21:00:37.882Z'
$ws.Range("B403").Value = 'HTTP'
$ws.Range("C403").Value = 'GET'
$ws.Range("D403").Value = '/user/EA4C7814?roomID=105'
$ws.Range("E403").Value = 404
$ws.Range("F403").Value = 'ERROR'
# G403 = "" (blank in source data; Excel Value="" clears the cell)
# H403 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("I403").Value = "'105"
$ws.Range("J403").Value = 'UID EA4C7814 not found in /user'
$ws.Range("K403").Value = '::ffff:172.28.219.204'

# Row 404
$ws.Range("A404").Value = '2025-11-28T21:00:38.339Z'
$ws.Range("B404").Value = 'HTTP'
$ws.Range("C404").Value = 'GET'
$ws.Range("D404").Value = '/uid-name/EA4C7814'
$ws.Range("E404").Value = 404
$ws.Range("F404").Value = 'ERROR'
# G404 = "" (blank in source data; Excel Value="" clears the cell)
# H404 = "" (blank in source data; Excel Value="" clears the cell)
# I404 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("J404").Value = 'UID EA4C7814 not found in /uid-name'
$ws.Range("K404").Value = '::ffff:172.28.219.204'

# Row 405
$ws.Range("A405").Value = '2025-11-28T21:00:42.719Z'
$ws.Range("B405").Value = 'HTTP'
$ws.Range("C405").Value = 'GET'
$ws.Range("D405").Value = '/user/EA4C7814?roomID=105'
$ws.Range("E405").Value = 404
$ws.Range("F405").Value = 'ERROR'
# G405 = "" (blank in source data; Excel Value="" clears the cell)
# H405 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("I405").Value = "'105"
$ws.Range("J405").Value = 'UID EA4C7814 not found in /user'
$ws.Range("K405").Value = '::ffff:172.28.219.204'

# Row 406
$ws.Range("A406").Value = '2025-11-28T21:00:43.001Z'
$ws.Range("B406").Value = 'HTTP'
$ws.Range("C406").Value = 'GET'
$ws.Range("D406").Value = '/uid-name/EA4C7814'
$ws.Range("E406").Value = 404
$ws.Range("F406").Value = 'ERROR'
# G406 = "" (blank in source data; Excel Value="" clears the cell)
# H406 = "" (blank in source data; Excel Value="" clears the cell)
# I406 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("J406").Value = 'UID EA4C7814 not found in /uid-name'
$ws.Range("K406").Value = '::ffff:172.28.219.204'

# Row 407
$ws.Range("A407").Value = '2025-11-28T21:00:56.442Z'
$ws.Range("B407").Value = 'HTTP'
$ws.Range("C407").Value = 'GET'
$ws.Range("D407").Value = '/stats?user=YAHYA&password=YAHYA'
$ws.Range("E407").Value = 200
$ws.Range("F407").Value = 'OK'
$ws.Range("G407").Value = 'YAHYA'
$ws.Range("H407").Value = '591EF2D4'
# I407 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("J407").Value = 'Stats returned for user'
$ws.Range("K407").Value = '::ffff:172.28.219.204'

# Row 408
$ws.Range("A408").Value = '2025-11-28T21:01:17.110Z'
$ws.Range("B408").Value = 'HTTP'
$ws.Range("C408").Value = 'POST'
$ws.Range("D408").Value = '/register'
$ws.Range("E408").Value = 200
$ws.Range("F408").Value = 'OK'
$ws.Range("G408").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H408").Value = 'E9956AF6'
$ws.Range("I408").Value = "'105"
$ws.Range("J408").Value = 'Registration successful for new UID'
$ws.Range("K408").Value = '::ffff:172.28.219.204'

# Row 409
$ws.Range("A409").Value = '2025-11-28T21:01:19.295Z'
$ws.Range("B409").Value = 'HTTP'
$ws.Range("C409").Value = 'GET'
$ws.Range("D409").Value = '/user/E9956AF6?roomID=105'
$ws.Range("E409").Value = 200
$ws.Range("F409").Value = 'OK'
$ws.Range("G409").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H409").Value = 'E9956AF6'
$ws.Range("I409").Value = "'105"
$ws.Range("J409").Value = 'Access granted and counter incremented'
$ws.Range("K409").Value = '::ffff:172.28.219.204'

# Row 410
$ws.Range("A410").Value = '2025-11-28T21:01:19.436Z'
$ws.Range("B410").Value = 'HTTP'
$ws.Range("C410").Value = 'GET'
$ws.Range("D410").Value = '/uid-name/E9956AF6'
$ws.Range("E410").Value = 200
$ws.Range("F410").Value = 'OK'
$ws.Range("G410").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H410").Value = 'E9956AF6'
# I410 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("J410").Value = 'UID to username lookup success'
$ws.Range("K410").Value = '::ffff:172.28.219.204'

# Row 411
$ws.Range("A411").Value = '2025-11-28T21:01:24.691Z'
$ws.Range("B411").Value = 'HTTP'
$ws.Range("C411").Value = 'GET'
$ws.Range("D411").Value = '/user/E9956AF6?roomID=105'
$ws.Range("E411").Value = 200
$ws.Range("F411").Value = 'OK'
$ws.Range("G411").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H411").Value = 'E9956AF6'
$ws.Range("I411").Value = "'105"
$ws.Range("J411").Value = 'Access granted and counter incremented'
$ws.Range("K411").Value = '::ffff:172.28.219.204'

# Row 412
$ws.Range("A412").Value = '2025-11-28T21:01:24.744Z'
$ws.Range("B412").Value = 'HTTP'
$ws.Range("C412").Value = 'GET'
$ws.Range("D412").Value = '/uid-name/E9956AF6'
$ws.Range("E412").Value = 200
$ws.Range("F412").Value = 'OK'
$ws.Range("G412").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H412").Value = 'E9956AF6'
# I412 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("J412").Value = 'UID to username lookup success'
$ws.Range("K412").Value = '::ffff:172.28.219.204'

# Row 413
$ws.Range("A413").Value = '2025-11-28T21:01:31.476Z'
$ws.Range("B413").Value = 'HTTP'
$ws.Range("C413").Value = 'GET'
$ws.Range("D413").Value = '/user/E9956AF6?roomID=105'
$ws.Range("E413").Value = 200
$ws.Range("F413").Value = 'OK'
$ws.Range("G413").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H413").Value = 'E9956AF6'
$ws.Range("I413").Value = "'105"
$ws.Range("J413").Value = 'Access granted and counter incremented'
$ws.Range("K413").Value = '::ffff:172.28.219.204'

# Row 414
$ws.Range("A414").Value = '2025-11-28T21:01:31.524Z'
$ws.Range("B414").Value = 'HTTP'
$ws.Range("C414").Value = 'GET'
$ws.Range("D414").Value = '/uid-name/E9956AF6'
$ws.Range("E414").Value = 200
$ws.Range("F414").Value = 'OK'
$ws.Range("G414").Value = 'hhhhhhhhhhhhhhhhhhhhhhh'
$ws.Range("H414").Value = 'E9956AF6'
# I414 = "" (blank in source data; Excel Value="" clears the cell)
$ws.Range("J414").Value = 'UID to username lookup success'
$ws.Range("K414").Value = '::ffff:172.28.219.204'

